$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -12.4354
$ws.Range("B3").Value = 6.285100000000007
$ws.Range("B14").Value = 5.674000000000003
$ws.Range("B21").Value = 9.643800000000004
$ws.Range("B23").Value = 9.152900000000006
$ws.Range("B25").Value = 5.448799999999999
$ws.Range("C25").Value = -13.3119
$ws.Range("B26").Value = 4.838800000000008
$ws.Range("C27").Value = -12.878
$ws.Range("B29").Value = 5.088300000000001
$ws.Range("C31").Value = -13.2368
$ws.Range("C39").Value = -12.56100000000001
$ws.Range("C48").Value = -11.24029999999999
$ws.Range("C51").Value = -11.2826
$ws.Range("C52").Value = -11.2416
$ws.Range("B53").Value = 5.0895
$ws.Range("C55").Value = -13.66119999999999
$ws.Range("C56").Value = -12.79729999999999
$ws.Range("B57").Value = 4.530799999999997
$ws.Range("C57").Value = -13.2919
$ws.Range("B59").Value = 4.748699999999997
$ws.Range("B69").Value = 5.666499999999996
$ws.Range("C73").Value = -13.06550000000001
$ws.Range("B79").Value = 8.928100000000002
$ws.Range("B83").Value = 5.241399999999999
$ws.Range("C89").Value = -10.6148
$ws.Range("C90").Value = -12.4586
$ws.Range("B91").Value = 4.967299999999999
$ws.Range("C92").Value = -11.0035
$ws.Range("B93").Value = 5.760600000000002
